$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 177; existing rows 177:203 shift down to 178:204.
$ws.Rows("177:177").Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

# Copy formatting (number format / style) from the row above it (old row 177, now 178)
# so the new row matches the rest of the table (e.g. date format in column D).
# Use the used-column range (A:R) rather than EntireRow so the sheet's used
# dimension doesn't get inflated to the full row width.
$ws.Range("A178:R178").Copy()
$ws.Range("A177:R177").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Populate the new row's data (matches the row for 2021-10-05 inserted into the series).
$ws.Range("A177").Value = 3
$ws.Range("B177").Value = "Femacal de La Calera"
$ws.Range("C177").Value = "Coquimbo"
$ws.Range("D177").Value = 44474
$ws.Range("E177").Value = 5
$ws.Range("F177").Value = 100112009
$ws.Range("G177").Value = "Acelga"
$ws.Range("H177").Value = "Sin especificar"
$ws.Range("I177").Value = "Primera"
$ws.Range("J177").Value = 240
$ws.Range("K177").Value = 2000
$ws.Range("L177").Value = 2200
$ws.Range("M177").Value = 2092
$ws.Range("N177").Value = "`$/docena de atados (6 kilos)"
$ws.Range("O177").Value = "Provincia de Quillota"
$ws.Range("P177").Value = 349
$ws.Range("Q177").Value = 6
$ws.Range("R177").Value = "Hortaliza"
